$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/29/2024  Through  2/4/2024"

# --- Data table updates (rows 14-30) ---
$ws.Range("L14").Value = -100
$ws.Range("L14").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("N15").Value = -87.5
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 2
$ws.Range("F16").Value = 20
$ws.Range("H16").Value = 150
$ws.Range("I16").Value = 23
$ws.Range("J16").Value = 10
$ws.Range("K16").Value = 130
$ws.Range("L16").Value = 53.333333333333
$ws.Range("M16").Value = 43.75
$ws.Range("N16").Value = -80.172413793103
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 100
$ws.Range("I17").Value = 17
$ws.Range("J17").Value = 11
$ws.Range("K17").Value = 54.545454545454
$ws.Range("L17").Value = 88.888888888888
$ws.Range("M17").Value = 6.25
$ws.Range("N17").Value = -58.536585365853
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 33.333333333333
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -9.090909090909
$ws.Range("I18").Value = 16
$ws.Range("J18").Value = 15
$ws.Range("K18").Value = 6.666666666666
$ws.Range("L18").Value = 14.285714285714
$ws.Range("M18").Value = -23.809523809523
$ws.Range("N18").Value = -92.694063926940
$ws.Range("C19").Value = 32
$ws.Range("D19").Value = 34
$ws.Range("E19").Value = -5.882352941176
$ws.Range("F19").Value = 148
$ws.Range("G19").Value = 148
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 182
$ws.Range("J19").Value = 190
$ws.Range("K19").Value = -4.210526315789
$ws.Range("L19").Value = 34.814814814814
$ws.Range("M19").Value = 13.043478260869
$ws.Range("N19").Value = -73.546511627907
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = -50
$ws.Range("I20").Value = 4
$ws.Range("J20").Value = 8
$ws.Range("K20").Value = -50
$ws.Range("L20").Value = -63.636363636363
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = -91.666666666666
$ws.Range("C21").Value = 45
$ws.Range("D21").Value = 45
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 196
$ws.Range("G21").Value = 180
$ws.Range("H21").Value = 8.888888888888
$ws.Range("I21").Value = 243
$ws.Range("J21").Value = 234
$ws.Range("K21").Value = 3.846153846153
$ws.Range("L21").Value = 30.645161290322
$ws.Range("M21").Value = 8.482142857142
$ws.Range("N21").Value = -78.380782918149
$ws.Range("C22").Value = 2
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("D22").Value = 2
$ws.Range("D22").NumberFormat = '#,##0'
$ws.Range("E22").Value = 0
$ws.Range("E22").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 6
$ws.Range("H22").Value = -16.666666666666
$ws.Range("I22").Value = 6
$ws.Range("J22").Value = 6
$ws.Range("L22").Value = 50
$ws.Range("D23").Value = 1
$ws.Range("D23").NumberFormat = '#,##0'
$ws.Range("E23").Value = -100
$ws.Range("E23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G23").Value = 1
$ws.Range("G23").NumberFormat = '#,##0'
$ws.Range("H23").Value = -100
$ws.Range("H23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J23").Value = 1
$ws.Range("J23").NumberFormat = '#,##0'
$ws.Range("K23").Value = -100
$ws.Range("K23").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C24").Value = 43
$ws.Range("D24").Value = 56
$ws.Range("E24").Value = -23.214285714285
$ws.Range("F24").Value = 220
$ws.Range("G24").Value = 230
$ws.Range("H24").Value = -4.347826086956
$ws.Range("I24").Value = 270
$ws.Range("J24").Value = 261
$ws.Range("K24").Value = 3.448275862068
$ws.Range("L24").Value = 83.673469387755
$ws.Range("M24").Value = 58.823529411764
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 47
$ws.Range("H25").Value = -12.765957446808
$ws.Range("I25").Value = 59
$ws.Range("J25").Value = 60
$ws.Range("K25").Value = -1.666666666666
$ws.Range("L25").Value = 55.263157894736
$ws.Range("M25").Value = 15.686274509803
$ws.Range("D26").Value = 1
$ws.Range("D26").NumberFormat = '#,##0'
$ws.Range("E26").Value = -100
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G26").Value = 1
$ws.Range("G26").NumberFormat = '#,##0'
$ws.Range("H26").Value = 0
$ws.Range("H26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("J26").Value = 1
$ws.Range("J26").NumberFormat = '#,##0'
$ws.Range("K26").Value = 0
$ws.Range("K26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("D27").Value = 5
$ws.Range("D27").NumberFormat = '#,##0'
$ws.Range("E27").Value = -80
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 3
$ws.Range("J27").Value = 6
$ws.Range("K27").Value = -50
$ws.Range("L27").Value = -66.666666666666
$ws.Range("D30").Value = 1
$ws.Range("D30").NumberFormat = '#,##0'
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("G30").Value = 2
$ws.Range("J30").Value = 2
